$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the "Totals" row (5-col summary row: [blank][blank]["Totals"][amount][blank])
# and insert a new grant row immediately above it, matching the new
# "co-PI, Lu, Y. (PI), Pyrcz, M. (co-PI)" / Hildebrand Seed Grant entry.
$totalsRowIndex = -1
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $row = $t.Rows.Item($i)
    if ($row.Cells.Item(3).Range.Text.TrimEnd([char]7, [char]13) -eq "Totals") {
        $totalsRowIndex = $i
        break
    }
}

$totalsRow = $t.Rows.Item($totalsRowIndex)
$newRow = $t.Rows.Add($totalsRow)

$newRow.Cells.Item(1).Range.Text = 'co-PI, Lu, Y. (PI), Pyrcz, M. (co-PI)'
$newRow.Cells.Item(2).Range.Text = 'Unconventional Well Optimization based on Machine Learning'
$newRow.Cells.Item(3).Range.Text = 'Hildebrand Seed Grant'
$newRow.Cells.Item(4).Range.Text = '$75,000 ($25,000)'
$newRow.Cells.Item(5).Range.Text = '9/2023-8/2025'

# Update the grand-totals amount to reflect the newly added grant.
$d.Content.Find.Execute('$14,019,591 ($3,780,350)', $true, $false, $false, $false, $false, `
                         $true, 1, $false, '$14,094,591 ($3,805,350)', 2)

# Update the "awarded in rank" subtotal the same way.
$d.Content.Find.Execute('$3,638,025 ($1,256,832)', $true, $false, $false, $false, $false, `
                         $true, 1, $false, '$3,713,025 ($1,281,832)', 2)
